$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (A:D) right by one column into (B:E), preserving the
# original column widths of B, C, D (which land on C, D, E) untouched.
$ws.Columns.Item(1).Insert()

# New column A width (closest achievable to the authored 54.552101 given
# this runtime's character-width quantization grid).
$ws.Columns.Item(1).ColumnWidth = 53.66

# Header row
$ws.Range("B1").Value = "Valid"
$ws.Range("C1").Value = "T"
$ws.Range("D1").Value = "Z"
$ws.Range("E1").Value = "p-value"

# Row labels (column A) and data (columns B:E)
$ws.Range("A2").Value = "CyclomaticComplexity(CC) & CyclomaticComplexity(CC)"
$ws.Range("B2").Value = 59
$ws.Range("C2").Value = 731.000000
$ws.Range("D2").Value = 1.162388
$ws.Range("E2").Value = 0.245079

$ws.Range("A3").Value = "CyclomaticComplexity(CC) & EffortToImplement"
$ws.Range("B3").Value = 75
$ws.Range("C3").Value = 1406.000000
$ws.Range("D3").Value = 0.100331
$ws.Range("E3").Value = 0.920082

$ws.Range("A4").Value = "MaintainabilityIndex & MaintainabilityIndex"
$ws.Range("B4").Value = 61
$ws.Range("C4").Value = 722.500000
$ws.Range("D4").Value = 1.601759
$ws.Range("E4").Value = 0.109210

$ws.Range("A5").Value = "NbUniqueOperands & NbUniqueOperands"
$ws.Range("B5").Value = 58
$ws.Range("C5").Value = 697.000000
$ws.Range("D5").Value = 1.227162
$ws.Range("E5").Value = 0.219763

$ws.Range("A6").Value = "NbOperands & NbOperands"
$ws.Range("B6").Value = 49
$ws.Range("C6").Value = 427.000000
$ws.Range("D6").Value = 1.845223
$ws.Range("E6").Value = 0.065006

$ws.Range("A7").Value = "NbOperands & EffortToImplement"
$ws.Range("B7").Value = 75
$ws.Range("C7").Value = 1100.000000
$ws.Range("D7").Value = 1.716181
$ws.Range("E7").Value = 0.086130

$ws.Range("A8").Value = "NbUniqueOperators & NbUniqueOperators"
$ws.Range("B8").Value = 63
$ws.Range("C8").Value = 958.000000
$ws.Range("D8").Value = 0.342305
$ws.Range("E8").Value = 0.732121

$ws.Range("A9").Value = "NbOperators & NbOperators"
$ws.Range("B9").Value = 43
$ws.Range("C9").Value = 458.000000
$ws.Range("D9").Value = 0.181124
$ws.Range("E9").Value = 0.856270

$ws.Range("A10").Value = "NbOperators & EffortToImplement"
$ws.Range("B10").Value = 75
$ws.Range("C10").Value = 1393.500000
$ws.Range("D10").Value = 0.166338
$ws.Range("E10").Value = 0.867891

$ws.Range("A11").Value = "ProgramLength & ProgramLength"
$ws.Range("B11").Value = 65
$ws.Range("C11").Value = 836.000000
$ws.Range("D11").Value = 1.545512
$ws.Range("E11").Value = 0.122223

$ws.Range("A12").Value = "VocabularySize & VocabularySize"
$ws.Range("B12").Value = 68
$ws.Range("C12").Value = 1100.000000
$ws.Range("D12").Value = 0.446056
$ws.Range("E12").Value = 0.655557

$ws.Range("A13").Value = "ProgramVolume & ProgramVolume"
$ws.Range("B13").Value = 73
$ws.Range("C13").Value = 1091.000000
$ws.Range("D13").Value = 1.426616
$ws.Range("E13").Value = 0.153692

$ws.Range("A14").Value = "DifficultyLevel & DifficultyLevel"
$ws.Range("B14").Value = 9
$ws.Range("C14").Value = 21.000000
$ws.Range("D14").Value = 0.177705
$ws.Range("E14").Value = 0.858955

$ws.Range("A15").Value = "ProgramLevel & ProgramLevel"
$ws.Range("B15").Value = 6
$ws.Range("C15").Value = 6.000000
$ws.Range("D15").Value = 0.943456
$ws.Range("E15").Value = 0.345448

$ws.Range("A16").Value = "EffortToImplement & CyclomaticComplexity(CC)"
$ws.Range("B16").Value = 75
$ws.Range("C16").Value = 1377.000000
$ws.Range("D16").Value = 0.253467
$ws.Range("E16").Value = 0.799908

$ws.Range("A17").Value = "EffortToImplement & NbOperands"
$ws.Range("B17").Value = 75
$ws.Range("C17").Value = 1115.000000
$ws.Range("D17").Value = 1.636973
$ws.Range("E17").Value = 0.101637

$ws.Range("A18").Value = "EffortToImplement & NbOperators"
$ws.Range("B18").Value = 75
$ws.Range("C18").Value = 1399.500000
$ws.Range("D18").Value = 0.134654
$ws.Range("E18").Value = 0.892885

$ws.Range("A19").Value = "EffortToImplement & EffortToImplement"
$ws.Range("B19").Value = 9
$ws.Range("C19").Value = 21.000000
$ws.Range("D19").Value = 0.177705
$ws.Range("E19").Value = 0.858955

$ws.Range("A20").Value = "TimeToImplement & TimeToImplement"
$ws.Range("B20").Value = 9
$ws.Range("C20").Value = 21.000000
$ws.Range("D20").Value = 0.177705
$ws.Range("E20").Value = 0.858955
